$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclaimer text date from 2021-05-24 to 2021-05-25
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-25 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.247252227644238
$ws.Range("E2").Value = 0.004295390715347791

$ws.Range("D3").Value = 0.4988608316239047
$ws.Range("E3").Value = -0.005760670332547679

$ws.Range("D4").Value = 0.09640845404457238
$ws.Range("E4").Value = 0.003332679866692834

$ws.Range("D5").Value = 0.1013337051391851
$ws.Range("E5").Value = -0.007334963325183574

$ws.Range("D6").Value = 0.05614478154809992
$ws.Range("E6").Value = -0.007034263671431784

$ws.Range("E7").Value = -0.002628645564045518
